$d = $word.ActiveDocument

# --- 1. Title replacement (appears twice: H1 heading and the bold blurb line) ---
$d.Content.Find.Execute(
    "Play Jumping Jokers Free - Review of Novomatic's Classic Fruit Slot",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Jumping Jokers Online for Free", 2)

# --- 2. Meta description (italic line) ---
$d.Content.Find.Execute(
    "Read our review of Jumping Jokers by Novomatic, a classic fruit slot game with Expanding Wild Symbols and Mystery Scatter Symbols. Play for free.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Jumping Jokers and play this classic fruit slot game for free.", 2)

# --- 3. "What we like" bullet list rewording ---
$d.Content.Find.Execute(
    "Expanding Wild Symbols increase the chances of winning",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Simple and classic fruit slot game", 2)

$d.Content.Find.Execute(
    "Mystery Scatter Symbol offers special prize",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Expanding wild symbols increase chances of winning", 2)

# Insert a new bullet paragraph ("Mystery scatter symbol offers special prizes")
# right after the paragraph that now reads "Expanding wild symbols increase chances of winning".
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Expanding wild symbols increase chances of winning")) {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "Mystery scatter symbol offers special prizes"
        break
    }
}

# Remove the old "Suitable for fans of traditional fruit slots" bullet entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Suitable for fans of traditional fruit slots")) {
        $p.Range.Delete()
        break
    }
}

# --- 4. "What we don't like" bullet list ---
$d.Content.Find.Execute(
    "No bonus features available",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lacks bonus features", 2)

# Remove the "Limited number of paylines" bullet entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Limited number of paylines")) {
        $p.Range.Delete()
        break
    }
}
